$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# B7 ("Experimental" row): set literal text "false" (not boolean FALSE).
# A direct Value = "false" assignment gets auto-coerced by Excel into a
# Boolean. Route it through a text formula, then freeze it back to a
# static value with PasteSpecial(xlPasteValues) so it lands as a plain
# shared-string cell (keeps the existing cell style, no quote-prefix).
$ws.Range("B7").Formula = '="false"'
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"
$ws.Range("B17").Value = "Supported wearable device vendors for data integration"
